$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-7
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45206
}
